$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge the "Et" and "Balık" product-group categories into a single
# "Et ve Balık" category for rows 10-20 (column C).
$ws.Range("C10:C20").Value = "Et ve Balık"

# Update the selection to match the newly edited range.
$ws.Range("C10:C20").Select()
